$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (2023-08-19)
$ws.Range("B2").Value = 0.3464964993005633
$ws.Range("C2").Value = 0.3375848360084654
$ws.Range("D2").Value = 3.082599426703578
$ws.Range("E2").Value = 0.4998867070740569
$ws.Range("G2").Value = 4.266567469086664

# Row 3 (2023-05-02)
$ws.Range("B3").Value = 3.182878228561681
$ws.Range("C3").Value = 1.65323645889881
$ws.Range("D3").Value = 3.082599426703578
$ws.Range("E3").Value = 0.4998867070740569
$ws.Range("G3").Value = 8.418600821238126
